$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the two new product rows below the existing "Samsung m32 128gb mobile" entry in A1.
$ws.Range("A2").Value = "apple mobile 14 pro max"
$ws.Range("A3").Value = "redmi note 12 pro plus 5g"

# Match the author's final selection/active cell (A3) recorded in the saved file.
$ws.Range("A3").Select()
